# Update "旅客运输平均运距" (average passenger transport distance) sheet:
#   - drop the oldest decade (2000年..2009年, old rows 2-11)
#   - this shifts the remaining 2010年..2020年 rows (old rows 12-22) up to
#     become the new rows 2-12
#   - append a new 2021年 row (row 13) with its own figures
#   - the new used range becomes A1:F13 (Excel updates this automatically)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2000年-2009年 rows entirely; everything below shifts up.
$ws.Rows("2:11").Delete()

# Copy the formatting (bold/border/alignment) used by the year-label column
# down onto the new last row before we populate it.
$ws.Cells.Item(12, 1).Copy()
$ws.Cells.Item(13, 1).PasteSpecial(-4122)  # xlPasteFormats

# Populate the new 2021年 row.
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 71.3110033242108
$ws.Cells.Item(13, 3).Value = 237.976408872398
$ws.Cells.Item(13, 4).Value = 1482.14078902427
$ws.Cells.Item(13, 5).Value = 20.2678229359027
$ws.Cells.Item(13, 6).Value = 366.343481182132
